$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$cells = $newRow.Cells
$cells.Item(1).Range.Text = "Entity"
$cells.Item(2).Range.Text = "Directory"
$cell3 = $cells.Item(3)
$cell3.Range.Text = "/models/directory.js"

$rngFind = $cell3.Range
$found = $rngFind.Find.Execute("directory", $true, $false, $false, $false, $false, $true, 1, $false)
$s = $rngFind.Start
$e = $rngFind.End
$rng2 = $d.Range($s, $e)
$rng2.LanguageIDFarEast = 2052
Write-Output "step2: $($d.Content.Text)"
